$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "('Overtaker', ['{1}{U}', 'Creature — Merfolk Spellshaper', '{3}{U}, {T}, Discard a card: Untap target creature and gain control of it until end of turn. That creature gains haste until end of turn.', '1/1'])"
$ws.Range("A3:A6").ClearContents()
